$wb = $excel.ActiveWorkbook

# --- "About" sheet: update the source/last-updated date in C1 ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45379

# --- "FPIEBP" sheet: update hard coal's production/import/export priorities ---
$wsFpiebp = $wb.Worksheets.Item("FPIEBP")
$wsFpiebp.Range("B3").Value = 1
$wsFpiebp.Range("C3").Value = 3
$wsFpiebp.Range("D3").Value = 2

# Move the active selection on the FPIEBP sheet to E3 (matches the saved cursor position)
$wsFpiebp.Activate()
$wsFpiebp.Range("E3").Select()
